$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before current row 2 (pushing Tomate..Pimenton down to rows 8-13)
$ws.Rows("2:7").Insert()

# New ingredient rows inserted at the top (rows 2-7)
$ws.Range("A2").Value = "Carne de Res"
$ws.Range("B2").Value = 8000
$ws.Range("C2").Value = 20

$ws.Range("A3").Value = "Carne de Pollo"
$ws.Range("B3").Value = 7000
$ws.Range("C3").Value = 20

$ws.Range("A4").Value = "Mixta"
$ws.Range("B4").Value = 10000
$ws.Range("C4").Value = 20

$ws.Range("A5").Value = "Pan Integral"
$ws.Range("B5").Value = 2000
$ws.Range("C5").Value = 20

$ws.Range("A6").Value = "Pan Blanco"
$ws.Range("B6").Value = 2000
$ws.Range("C6").Value = 20

$ws.Range("A7").Value = "Pan de Oregano"
$ws.Range("B7").Value = 2000
$ws.Range("C7").Value = 20

# New ingredient rows appended at the bottom (rows 14-17)
$ws.Range("A14").Value = "Mayonesa"
$ws.Range("B14").Value = 200
$ws.Range("C14").Value = 20

$ws.Range("A15").Value = "Ketchup"
$ws.Range("B15").Value = 200
$ws.Range("C15").Value = 20

$ws.Range("A16").Value = "Piña"
$ws.Range("B16").Value = 200
$ws.Range("C16").Value = 20

$ws.Range("A17").Value = "Cebolla Dulce"
$ws.Range("B17").Value = 200
$ws.Range("C17").Value = 20

# Update selection to match final state
$ws.Range("F9").Select()
